$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: move the existing "AVG" column (currently G) to the new column I ---
# Copy header text.
$ws.Cells.Item(1, 9).HorizontalAlignment = -4108
$ws.Cells.Item(1, 9).Value2 = $ws.Cells.Item(1, 7).Value2

# Copy the AVERAGE(...) formulas (and their cached results) row by row.
for ($r = 2; $r -le 15; $r++) {
    $srcCell = $ws.Cells.Item($r, 7)
    $dstCell = $ws.Cells.Item($r, 9)
    $dstCell.HorizontalAlignment = -4108
    $dstCell.Formula = "=AVERAGE(B$r`:F$r)"
}

# --- Step 2: clear the old column G content (it will be replaced by SD5009 below) ---
$ws.Range("G1:G15").ClearContents()

# --- Step 3: write the new SD5009 (col G) and SD5010 (col H) survey columns ---
$ws.Cells.Item(1, 7).HorizontalAlignment = -4108
$ws.Cells.Item(1, 7).Value2 = "SD5009"
$ws.Cells.Item(1, 8).HorizontalAlignment = -4108
$ws.Cells.Item(1, 8).Value2 = "SD5010"

$sd5009 = @(0, 100, 92, 86, 84, 40, 0)
for ($i = 0; $i -lt $sd5009.Length; $i++) {
    $r = 2 + $i
    $c = $ws.Cells.Item($r, 7)
    $c.HorizontalAlignment = -4108
    $c.Value2 = $sd5009[$i]
}

$sd5010 = @(0, 2, 10, 30, 74, 100, 94, 84, 66, 34, 0)
for ($i = 0; $i -lt $sd5010.Length; $i++) {
    $r = 2 + $i
    $c = $ws.Cells.Item($r, 8)
    $c.HorizontalAlignment = -4108
    $c.Value2 = $sd5010[$i]
}

$wb.Save()
